# Update cryptocurrency Price (column D) and Volume(1h) (column E) cells
# for rows 2-51 on the active sheet to reflect the latest scraped values.
# Numeric-looking Price values are prefixed with a leading apostrophe so Excel
# stores them as literal text (matching the workbook's inlineStr cells) instead
# of auto-converting them to numbers (which would also drop significant digits).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.179.82"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "1.868.51"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'306.15"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("D7").Value = "'0.5174"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").Value = "'0.3746"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").Value = "'0.07155"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("D10").Value = "'0.8920"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").Value = "'20.77"
$ws.Range("E11").Value = "  -2.16%  "
$ws.Range("D12").Value = "1.874.19"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "'0.07538"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "'5.306"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "'89.76"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("D16").Value = "'0.9998"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").Value = "'14.13"
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("D19").Value = "'0.9993"
$ws.Range("D20").Value = "27.194.65"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "'5.003"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "2.104.81"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("E23").Value = "  -3.39%  "
$ws.Range("D24").Value = "'6.473"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").Value = "'1.837"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "'146.29"
$ws.Range("D27").Value = "'17.98"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").Value = "'2.090"
$ws.Range("E28").Value = "  -3.12%  "
$ws.Range("D29").Value = "'112.98"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "'4.665"
$ws.Range("E30").Value = "  -3.69%  "
$ws.Range("D31").Value = "'4.698"
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("D32").Value = "'0.09261"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("D33").Value = "'0.05130"
$ws.Range("E33").Value = "  -3.15%  "
$ws.Range("D34").Value = "'3.082"
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("D35").Value = "'1.163"
$ws.Range("E35").Value = "  -5.16%  "
$ws.Range("D36").Value = "'0.7277"
$ws.Range("E36").Value = "  -6.80%  "
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("D38").Value = "'3.116"
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("D39").Value = "'2.511"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").Value = "'0.5312"
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").Value = "'6.526"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").Value = "'116.16"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'8.323"
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").Value = "'0.9991"
$ws.Range("D47").Value = "'0.4626"
$ws.Range("E47").Value = "  -4.10%  "
$ws.Range("D48").Value = "'9.986"
$ws.Range("E48").Value = "  -4.70%  "
$ws.Range("D49").Value = "'1.563"
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "'36.75"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "'63.74"
$ws.Range("E51").Value = "  -4.49%  "
